# fix(gui) step 1 and 2
# - Advance the sheet's date stamp in A1 by one day.
# - Update "CON TOPE" and "SIN TOPE" price lists (column D) to the new prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: bump the date in A1
$ws.Range("A1").Value = 45309

# Step 2: update prices in column D
$ws.Range("D34").Value = 1996.418
$ws.Range("D35").Value = 2661.888
$ws.Range("D36").Value = 3855.82
$ws.Range("D37").Value = 3279.202

$ws.Range("D41").Value = 2564.02
$ws.Range("D42").Value = 3387.579
$ws.Range("D43").Value = 4266.839
$ws.Range("D44").Value = 3729.338
